# Generate Report for handback
# Adds a new handback record (f980c316-8443-4d72-abec-f2cf5702fc64) as row 4
# to the Overview, zh-cn and de-de sheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$hyperlinkUnderline = 2
$hyperlinkColor = 15570276  # BGR-encoded RGB(0x64,0x95,0xED) == FF6495ED

function Style-AsHyperlink($range) {
    $range.Font.Underline = $hyperlinkUnderline
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A4").Value = "f980c316-8443-4d72-abec-f2cf5702fc64.md"
$ws1.Range("B4").Value = "Handed back: in sync with en-US"
$ws1.Range("C4").Value = "Handed back: in sync with en-US"

$ws1.Hyperlinks.Add(
    $ws1.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f980c316-handback/e2e/f980c316-8443-4d72-abec-f2cf5702fc64.md",
    "",
    "",
    "f980c316-8443-4d72-abec-f2cf5702fc64.md"
) | Out-Null
Style-AsHyperlink $ws1.Range("A4")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A4").Value = "f980c316-8443-4d72-abec-f2cf5702fc64.md"
$ws2.Range("B4").Value = "Handed back: in sync with en-US"
$ws2.Range("C4").Value = "f980c316-8443-4d72-abec-f2cf5702fc64.b72f8c30b1d118c6125851d67adbe0a3f38c2e2a.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-01-28 09:01:12"
$ws2.Range("E4").Value = "f980c316-8443-4d72-abec-f2cf5702fc64.md"
$ws2.Range("F4").Value = "f980c316-8443-4d72-abec-f2cf5702fc64.b72f8c30b1d118c6125851d67adbe0a3f38c2e2a.zh-cn.xlf"
$ws2.Range("G4").Value = "2016-01-28 09:02:05"
$ws2.Range("H4").Value = "Include"

$ws2.Hyperlinks.Add(
    $ws2.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f980c316-handback/e2e/f980c316-8443-4d72-abec-f2cf5702fc64.md",
    "",
    "",
    "f980c316-8443-4d72-abec-f2cf5702fc64.md"
) | Out-Null
Style-AsHyperlink $ws2.Range("A4")

$ws2.Hyperlinks.Add(
    $ws2.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f980c316-handoff/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f980c316-8443-4d72-abec-f2cf5702fc64.b72f8c30b1d118c6125851d67adbe0a3f38c2e2a.zh-cn.xlf",
    "",
    "",
    "f980c316-8443-4d72-abec-f2cf5702fc64.b72f8c30b1d118c6125851d67adbe0a3f38c2e2a.zh-cn.xlf"
) | Out-Null
Style-AsHyperlink $ws2.Range("C4")

$ws2.Hyperlinks.Add(
    $ws2.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f980c316-md/e2e/f980c316-8443-4d72-abec-f2cf5702fc64.md",
    "",
    "",
    "f980c316-8443-4d72-abec-f2cf5702fc64.md"
) | Out-Null
Style-AsHyperlink $ws2.Range("E4")

$ws2.Hyperlinks.Add(
    $ws2.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f980c316-handback/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f980c316-8443-4d72-abec-f2cf5702fc64.b72f8c30b1d118c6125851d67adbe0a3f38c2e2a.zh-cn.xlf",
    "",
    "",
    "f980c316-8443-4d72-abec-f2cf5702fc64.b72f8c30b1d118c6125851d67adbe0a3f38c2e2a.zh-cn.xlf"
) | Out-Null
Style-AsHyperlink $ws2.Range("F4")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A4").Value = "f980c316-8443-4d72-abec-f2cf5702fc64.md"
$ws3.Range("B4").Value = "Handed back: in sync with en-US"
$ws3.Range("C4").Value = "f980c316-8443-4d72-abec-f2cf5702fc64.b72f8c30b1d118c6125851d67adbe0a3f38c2e2a.de-de.xlf"
$ws3.Range("D4").Value = "2016-01-28 09:01:28"
$ws3.Range("E4").Value = "f980c316-8443-4d72-abec-f2cf5702fc64.md"
$ws3.Range("F4").Value = "f980c316-8443-4d72-abec-f2cf5702fc64.b72f8c30b1d118c6125851d67adbe0a3f38c2e2a.de-de.xlf"
$ws3.Range("G4").Value = "2016-01-28 09:02:26"
$ws3.Range("H4").Value = "Include"

$ws3.Hyperlinks.Add(
    $ws3.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f980c316-handback/e2e/f980c316-8443-4d72-abec-f2cf5702fc64.md",
    "",
    "",
    "f980c316-8443-4d72-abec-f2cf5702fc64.md"
) | Out-Null
Style-AsHyperlink $ws3.Range("A4")

$ws3.Hyperlinks.Add(
    $ws3.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f980c316-handoff/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f980c316-8443-4d72-abec-f2cf5702fc64.b72f8c30b1d118c6125851d67adbe0a3f38c2e2a.de-de.xlf",
    "",
    "",
    "f980c316-8443-4d72-abec-f2cf5702fc64.b72f8c30b1d118c6125851d67adbe0a3f38c2e2a.de-de.xlf"
) | Out-Null
Style-AsHyperlink $ws3.Range("C4")

$ws3.Hyperlinks.Add(
    $ws3.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/f980c316-md/e2e/f980c316-8443-4d72-abec-f2cf5702fc64.md",
    "",
    "",
    "f980c316-8443-4d72-abec-f2cf5702fc64.md"
) | Out-Null
Style-AsHyperlink $ws3.Range("E4")

$ws3.Hyperlinks.Add(
    $ws3.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f980c316-handback/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f980c316-8443-4d72-abec-f2cf5702fc64.b72f8c30b1d118c6125851d67adbe0a3f38c2e2a.de-de.xlf",
    "",
    "",
    "f980c316-8443-4d72-abec-f2cf5702fc64.b72f8c30b1d118c6125851d67adbe0a3f38c2e2a.de-de.xlf"
) | Out-Null
Style-AsHyperlink $ws3.Range("F4")

Write-Host "Handback report row added for f980c316-8443-4d72-abec-f2cf5702fc64"
